# Auto-generated edit script applying the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds text-formatted numbers (grouping dots, fixed decimals).
# Pre-format every Price cell we touch as Text so Excel does not silently
# reinterpret the literal string as a Number (which would drop separators /
# significant trailing zeros, e.g. "1.00" -> 1, "441.10" -> 441.1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.901.76"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "3.148.53"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "588.02"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "146.13"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.138.49"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +7.04%  "
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "37.07"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D16").Value = "3.649.21"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "63.608.31"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "3.131.28"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "465.92"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "7.46"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "13.04"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").Value = "81.39"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "9.11"
$ws.Range("E28").Value = "  +6.62%  "
$ws.Range("D29").Value = "2.69"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "7.02"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").Value = "27.02"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.0₃0868"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("E38").Value = "  -4.30%  "
$ws.Range("D39").Value = "6.01"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "50.39"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").Value = "441.10"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").Value = "8.72"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "2.916.43"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "0.274"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").Value = "36.22"
$ws.Range("D48").Value = "125.70"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").Value = "24.36"
$ws.Range("E51").Value = "  -1.05%  "
